$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D cells below contain price text that looks numeric (e.g. "306.76").
# Force Text number format first so Excel keeps them as literal strings
# (matching the source data which stores prices/volumes as text), rather
# than auto-converting to floating point numbers.
$textCells = @("D5","D7","D8","D9","D10","D11","D12","D14","D15","D16","D17","D18","D19","D21","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.194.33"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.904.42"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "306.76"
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "0.5277"
$ws.Range("E7").Value = "  +2.51%  "
$ws.Range("D8").Value = "0.3780"
$ws.Range("E8").Value = "  +1.59%  "
$ws.Range("D9").Value = "0.07251"
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("D10").Value = "21.17"
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("D11").Value = "0.9003"
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "0.08399"
$ws.Range("E12").Value = "  +10.16%  "
$ws.Range("D13").Value = "1.916.91"
$ws.Range("E13").Value = "  +1.49%  "
$ws.Range("D14").Value = "94.99"
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("D15").Value = "5.274"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "0.000008618"
$ws.Range("E17").Value = "  +1.27%  "
$ws.Range("D18").Value = "14.58"
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("D19").Value = "1.0000"
$ws.Range("D20").Value = "27.234.33"
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("D21").Value = "5.065"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").Value = "2.147.16"
$ws.Range("E22").Value = "  +1.79%  "
$ws.Range("D23").Value = "10.60"
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("D24").Value = "6.443"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").Value = "147.22"
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("D26").Value = "2.280"
$ws.Range("E26").Value = "  +5.76%  "
$ws.Range("D27").Value = "1.751"
$ws.Range("E27").Value = "  -2.35%  "
$ws.Range("D28").Value = "18.18"
$ws.Range("E28").Value = "  +0.78%  "
$ws.Range("D29").Value = "114.91"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "4.932"
$ws.Range("E30").Value = "  -1.84%  "
$ws.Range("D31").Value = "4.828"
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("D32").Value = "0.09296"
$ws.Range("E32").Value = "  +0.90%  "
$ws.Range("D33").Value = "0.8111"
$ws.Range("E33").Value = "  +6.96%  "
$ws.Range("D34").Value = "0.05069"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").Value = "1.239"
$ws.Range("E35").Value = "  +3.30%  "
$ws.Range("D36").Value = "2.953"
$ws.Range("E36").Value = "  -2.27%  "
$ws.Range("D37").Value = "3.366"
$ws.Range("E37").Value = "  +2.72%  "
$ws.Range("D38").Value = "2.625"
$ws.Range("E38").Value = "  +2.60%  "
$ws.Range("D39").Value = "0.5733"
$ws.Range("E39").Value = "  +1.42%  "
$ws.Range("D40").Value = "0.01988"
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("D41").Value = "1.075"
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("D42").Value = "6.646"
$ws.Range("E42").Value = "  +1.03%  "
$ws.Range("D43").Value = "8.986"
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("D44").Value = "117.53"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("D45").Value = "0.1515"
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("D46").Value = "0.4851"
$ws.Range("E46").Value = "  +1.02%  "
$ws.Range("D47").Value = "10.21"
$ws.Range("E47").Value = "  +0.64%  "
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("D49").Value = "1.618"
$ws.Range("E49").Value = "  +2.58%  "
$ws.Range("D50").Value = "37.48"
$ws.Range("E50").Value = "  +0.78%  "
$ws.Range("D51").Value = "63.73"
$ws.Range("E51").Value = "  +0.29%  "
